$wb = $excel.ActiveWorkbook

function Fill-Sheet($ws, $data) {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $row = $data[$i]
        $r = $i + 2
        for ($c = 1; $c -le 7; $c++) {
            $ws.Cells.Item($r, $c).Value = $row[$c - 1]
        }
    }
}

$data1 = @(
    ,@("Frac. Fácil ELO p3 c3", "0.504 (N = 180)", "0.021 (N = 210)", "0.097 (N = 140)", "0.063 (N = 354)", "0.975 (N = 208)", "0.05 (N = 1092)")
    ,@("Frac. Fácil ELO p4 c3", "0.282 (N = 180)", "0.028 (N = 210)", "0.028 (N = 140)", "0.135 (N = 354)", "0.902 (N = 208)", "0.009 (N = 1092)")
    ,@("Frac. Fácil ELO p5 c3", "0.03 (N = 180)", "0.161 (N = 210)", "0.089 (N = 140)", "0.22 (N = 354)", "0.953 (N = 208)", "0.025 (N = 1092)")
    ,@("Frac. Fácil ELO p6 c3", "0.031 (N = 180)", "0.064 (N = 210)", "0.111 (N = 140)", "0.318 (N = 354)", "0.866 (N = 208)", "0.019 (N = 1092)")
    ,@("Frac. Fácil ELO p3 c4", "0.346 (N = 150)", "0.428 (N = 180)", "0.339 (N = 120)", "0.458 (N = 304)", "0.85 (N = 178)", "0.331 (N = 932)")
    ,@("Frac. Fácil ELO p4 c4", "0.326 (N = 150)", "0.094 (N = 180)", "0.326 (N = 120)", "0.576 (N = 304)", "0.351 (N = 178)", "0.07 (N = 932)")
    ,@("Frac. Fácil ELO p5 c4", "0.009 (N = 150)", "0.342 (N = 180)", "0.559 (N = 120)", "0.519 (N = 304)", "0.374 (N = 178)", "0.027 (N = 932)")
    ,@("Frac. Fácil ELO p6 c4", "0.035 (N = 150)", "0.092 (N = 180)", "0.87 (N = 120)", "0.728 (N = 304)", "0.458 (N = 178)", "0.041 (N = 932)")
    ,@("Frac. Fácil ELO p3 c5", "0.45 (N = 120)", "0.153 (N = 150)", "0.938 (N = 100)", "0.322 (N = 254)", "0.434 (N = 148)", "0.32 (N = 772)")
    ,@("Frac. Fácil ELO p4 c5", "0.58 (N = 120)", "0.049 (N = 150)", "0.855 (N = 100)", "0.496 (N = 254)", "0.235 (N = 148)", "0.191 (N = 772)")
    ,@("Frac. Fácil ELO p5 c5", "0.126 (N = 120)", "0.061 (N = 150)", "0.594 (N = 100)", "0.582 (N = 254)", "0.195 (N = 148)", "0.051 (N = 772)")
    ,@("Frac. Fácil ELO p6 c5", "0.218 (N = 120)", "0.009 (N = 150)", "0.708 (N = 100)", "0.346 (N = 254)", "0.409 (N = 148)", "0.016 (N = 772)")
    ,@("Frac. Fácil ELO p3 c6", "0.152 (N = 90)", "0.435 (N = 120)", "0.92 (N = 80)", "0.423 (N = 204)", "0.603 (N = 118)", "0.513 (N = 612)")
    ,@("Frac. Fácil ELO p4 c6", "0.238 (N = 90)", "0.103 (N = 120)", "0.791 (N = 80)", "0.357 (N = 204)", "0.454 (N = 118)", "0.146 (N = 612)")
    ,@("Frac. Fácil ELO p5 c6", "0.419 (N = 90)", "0.064 (N = 120)", "0.683 (N = 80)", "0.234 (N = 204)", "0.352 (N = 118)", "0.121 (N = 612)")
    ,@("Frac. Fácil ELO p6 c6", "0.597 (N = 90)", "0.006 (N = 120)", "0.702 (N = 80)", "0.083 (N = 204)", "0.265 (N = 118)", "0.02 (N = 612)")
)

$data2 = @(
    ,@("Frac. Difícil ELO p3 c3", "0.996 (N = 180)", "0.31 (N = 210)", "0.633 (N = 140)", "0.526 (N = 354)", "0.615 (N = 208)", "0.863 (N = 1092)")
    ,@("Frac. Difícil ELO p4 c3", "0.998 (N = 180)", "0.207 (N = 210)", "0.472 (N = 140)", "0.926 (N = 354)", "0.159 (N = 208)", "0.834 (N = 1092)")
    ,@("Frac. Difícil ELO p5 c3", "0.933 (N = 180)", "0.03 (N = 210)", "0.596 (N = 140)", "0.933 (N = 354)", "0.044 (N = 208)", "0.356 (N = 1092)")
    ,@("Frac. Difícil ELO p6 c3", "0.806 (N = 180)", "0.252 (N = 210)", "0.379 (N = 140)", "0.74 (N = 354)", "0.34 (N = 208)", "0.537 (N = 1092)")
    ,@("Frac. Difícil ELO p3 c4", "0.166 (N = 150)", "0.206 (N = 180)", "0.051 (N = 120)", "0.894 (N = 304)", "0.74 (N = 178)", "0.444 (N = 932)")
    ,@("Frac. Difícil ELO p4 c4", "0.757 (N = 150)", "0.023 (N = 180)", "0.015 (N = 120)", "0.991 (N = 304)", "0.215 (N = 178)", "0.352 (N = 932)")
    ,@("Frac. Difícil ELO p5 c4", "0.192 (N = 150)", "0.013 (N = 180)", "0.162 (N = 120)", "0.939 (N = 304)", "0.091 (N = 178)", "0.073 (N = 932)")
    ,@("Frac. Difícil ELO p6 c4", "0.408 (N = 150)", "0.068 (N = 180)", "0.082 (N = 120)", "0.715 (N = 304)", "0.419 (N = 178)", "0.182 (N = 932)")
    ,@("Frac. Difícil ELO p3 c5", "0.057 (N = 120)", "0.117 (N = 150)", "0.273 (N = 100)", "0.576 (N = 254)", "0.631 (N = 148)", "0.249 (N = 772)")
    ,@("Frac. Difícil ELO p4 c5", "0.403 (N = 120)", "0.004 (N = 150)", "0.289 (N = 100)", "0.633 (N = 254)", "0.166 (N = 148)", "0.071 (N = 772)")
    ,@("Frac. Difícil ELO p5 c5", "0.072 (N = 120)", "0.025 (N = 150)", "0.596 (N = 100)", "0.755 (N = 254)", "0.103 (N = 148)", "0.088 (N = 772)")
    ,@("Frac. Difícil ELO p6 c5", "0.379 (N = 120)", "0.12 (N = 150)", "0.463 (N = 100)", "0.609 (N = 254)", "0.074 (N = 148)", "0.206 (N = 772)")
    ,@("Frac. Difícil ELO p3 c6", "0.127 (N = 90)", "0.276 (N = 120)", "0.554 (N = 80)", "0.347 (N = 204)", "0.375 (N = 118)", "0.284 (N = 612)")
    ,@("Frac. Difícil ELO p4 c6", "0.467 (N = 90)", "0.008 (N = 120)", "0.725 (N = 80)", "0.486 (N = 204)", "0.299 (N = 118)", "0.17 (N = 612)")
    ,@("Frac. Difícil ELO p5 c6", "0.456 (N = 90)", "0.01 (N = 120)", "0.942 (N = 80)", "0.683 (N = 204)", "0.227 (N = 118)", "0.242 (N = 612)")
    ,@("Frac. Difícil ELO p6 c6", "0.695 (N = 90)", "0.023 (N = 120)", "0.964 (N = 80)", "0.185 (N = 204)", "0.173 (N = 118)", "0.136 (N = 612)")
)

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

Fill-Sheet $ws1 $data1
Fill-Sheet $ws2 $data2

Write-Output ("sheet1 dim: " + $ws1.UsedRange.Address())
Write-Output ("sheet2 dim: " + $ws2.UsedRange.Address())
